# Recompute the cross-correlation table (macro_corr_byear_gr) using
# Newey-West standard errors: refresh the correlation coefficients and
# their significance stars in B2:M13.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table stores each statistic as text (e.g. "0.1", "-0.0", "0.21*")
# rather than a number, so trailing zeros/asterisks are preserved.
# Mark the data range as Text *before* writing so Excel does not
# auto-convert the numeric-looking strings into real numbers.
$ws.Range("B2:M13").NumberFormat = "@"

# Row 2
$ws.Range("B2").Value = '-0.08'
$ws.Range("C2").Value = '-0.03'
$ws.Range("D2").Value = '0.16'
$ws.Range("E2").Value = '0.21*'
$ws.Range("F2").Value = '0.23*'
$ws.Range("G2").Value = '0.05'
$ws.Range("H2").Value = '0.0'
$ws.Range("I2").Value = '0.19'
$ws.Range("J2").Value = '-0.01'
$ws.Range("K2").Value = '-0.04'
$ws.Range("L2").Value = '0.25**'
$ws.Range("M2").Value = '-0.01'

# Row 3
$ws.Range("B3").Value = '-0.16'
$ws.Range("C3").Value = '0.1'
$ws.Range("D3").Value = '0.18'
$ws.Range("E3").Value = '0.16'
$ws.Range("F3").Value = '0.16'
$ws.Range("G3").Value = '0.18'
$ws.Range("H3").Value = '0.06'
$ws.Range("I3").Value = '0.12'
$ws.Range("J3").Value = '-0.03'
$ws.Range("K3").Value = '-0.09'
$ws.Range("L3").Value = '0.17'
$ws.Range("M3").Value = '-0.01'

# Row 4
$ws.Range("B4").Value = '-0.14'
$ws.Range("C4").Value = '0.05'
$ws.Range("D4").Value = '0.21*'
$ws.Range("E4").Value = '0.17'
$ws.Range("F4").Value = '0.15'
$ws.Range("G4").Value = '0.26**'
$ws.Range("H4").Value = '0.1'
$ws.Range("I4").Value = '0.23*'
$ws.Range("J4").Value = '-0.04'
$ws.Range("K4").Value = '-0.2'
$ws.Range("L4").Value = '0.07'
$ws.Range("M4").Value = '-0.08'

# Row 5
$ws.Range("B5").Value = '-0.18'
$ws.Range("C5").Value = '-0.04'
$ws.Range("D5").Value = '0.2'
$ws.Range("E5").Value = '0.23*'
$ws.Range("F5").Value = '0.15'
$ws.Range("G5").Value = '0.17'
$ws.Range("H5").Value = '0.11'
$ws.Range("I5").Value = '0.18'
$ws.Range("J5").Value = '-0.07'
$ws.Range("K5").Value = '-0.12'
$ws.Range("L5").Value = '0.07'
$ws.Range("M5").Value = '-0.12'

# Row 6
$ws.Range("B6").Value = '-0.13'
$ws.Range("C6").Value = '-0.14'
$ws.Range("D6").Value = '0.05'
$ws.Range("E6").Value = '0.23*'
$ws.Range("F6").Value = '0.12'
$ws.Range("G6").Value = '-0.01'
$ws.Range("H6").Value = '-0.03'
$ws.Range("I6").Value = '0.14'
$ws.Range("J6").Value = '-0.08'
$ws.Range("K6").Value = '0.0'
$ws.Range("L6").Value = '0.12'
$ws.Range("M6").Value = '-0.18'

# Row 7
$ws.Range("B7").Value = '-0.05'
$ws.Range("C7").Value = '-0.14'
$ws.Range("D7").Value = '0.06'
$ws.Range("E7").Value = '0.22*'
$ws.Range("F7").Value = '0.23*'
$ws.Range("G7").Value = '-0.15'
$ws.Range("H7").Value = '-0.01'
$ws.Range("I7").Value = '0.09'
$ws.Range("J7").Value = '-0.02'
$ws.Range("K7").Value = '-0.03'
$ws.Range("L7").Value = '0.09'
$ws.Range("M7").Value = '-0.13'

# Row 8
$ws.Range("B8").Value = '0.03'
$ws.Range("C8").Value = '-0.15'
$ws.Range("D8").Value = '-0.07'
$ws.Range("E8").Value = '0.18'
$ws.Range("F8").Value = '0.19'
$ws.Range("G8").Value = '-0.03'
$ws.Range("H8").Value = '-0.1'
$ws.Range("I8").Value = '0.06'
$ws.Range("J8").Value = '-0.12'
$ws.Range("K8").Value = '-0.01'
$ws.Range("L8").Value = '0.02'
$ws.Range("M8").Value = '-0.2'

# Row 9
$ws.Range("B9").Value = '0.09'
$ws.Range("C9").Value = '-0.18'
$ws.Range("D9").Value = '-0.13'
$ws.Range("E9").Value = '0.07'
$ws.Range("F9").Value = '0.19'
$ws.Range("G9").Value = '-0.14'
$ws.Range("H9").Value = '-0.13'
$ws.Range("I9").Value = '-0.11'
$ws.Range("J9").Value = '-0.15'
$ws.Range("K9").Value = '0.05'
$ws.Range("L9").Value = '-0.11'
$ws.Range("M9").Value = '-0.26**'

# Row 10
$ws.Range("B10").Value = '-0.0'
$ws.Range("C10").Value = '-0.24*'
$ws.Range("D10").Value = '-0.09'
$ws.Range("E10").Value = '0.07'
$ws.Range("F10").Value = '0.11'
$ws.Range("G10").Value = '-0.14'
$ws.Range("H10").Value = '-0.1'
$ws.Range("I10").Value = '-0.03'
$ws.Range("J10").Value = '-0.01'
$ws.Range("K10").Value = '0.05'
$ws.Range("L10").Value = '-0.1'
$ws.Range("M10").Value = '-0.33***'

# Row 11
$ws.Range("B11").Value = '-0.13'
$ws.Range("C11").Value = '-0.23*'
$ws.Range("D11").Value = '-0.11'
$ws.Range("E11").Value = '-0.03'
$ws.Range("F11").Value = '0.01'
$ws.Range("G11").Value = '-0.12'
$ws.Range("H11").Value = '-0.12'
$ws.Range("I11").Value = '-0.1'
$ws.Range("J11").Value = '-0.0'
$ws.Range("K11").Value = '0.01'
$ws.Range("L11").Value = '-0.05'
$ws.Range("M11").Value = '-0.27**'

# Row 12
$ws.Range("B12").Value = '-0.15'
$ws.Range("C12").Value = '-0.19'
$ws.Range("D12").Value = '-0.2'
$ws.Range("E12").Value = '-0.11'
$ws.Range("F12").Value = '-0.09'
$ws.Range("G12").Value = '-0.05'
$ws.Range("H12").Value = '-0.12'
$ws.Range("I12").Value = '-0.08'
$ws.Range("J12").Value = '-0.0'
$ws.Range("K12").Value = '-0.05'
$ws.Range("L12").Value = '-0.12'
$ws.Range("M12").Value = '-0.34***'

# Row 13
$ws.Range("B13").Value = '-0.08'
$ws.Range("C13").Value = '-0.07'
$ws.Range("D13").Value = '-0.19'
$ws.Range("E13").Value = '-0.02'
$ws.Range("F13").Value = '-0.04'
$ws.Range("G13").Value = '0.02'
$ws.Range("H13").Value = '-0.08'
$ws.Range("I13").Value = '-0.06'
$ws.Range("J13").Value = '-0.01'
$ws.Range("K13").Value = '-0.1'
$ws.Range("L13").Value = '-0.09'
$ws.Range("M13").Value = '-0.3**'
